# Generate Report for Handoff
# Updates the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
# values for the most recently processed file (ba6e6927-...) on each sheet.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: column G = "Latest HO Xliff Generate Date", row 7 = ba6e6927 file
$wsOverview.Range("G7").Value = "2016-08-16 14:42:38"

# zh-cn sheet: column H = "Latest Handoff Datetime", row 7 = ba6e6927 file
$wsZhCn.Range("H7").Value = "2016-08-16 14:42:33"

# de-de sheet: column H = "Latest Handoff Datetime", row 7 = ba6e6927 file
$wsDeDe.Range("H7").Value = "2016-08-16 14:42:38"
